# Update benchmark: 2025-10-14 06:39:20 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BENCHMARK")

# Clear the YKB (column F) duplicate figures that mirrored column E/others.
$ws.Range("F3").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("F8").Value = ""
$ws.Range("F9").Value = ""
$ws.Range("F10").Value = ""

# Row 13 (GELEN SWIFT) updates.
$ws.Range("E13").Value = "Hesaba: Asgari 1 TL | Azami 8.700 TL"
$ws.Range("F13").Value = ""
$ws.Range("K13").Value = "Hesaba: Asgari 1 TL | Azami 53,19 TL"

# Row 14 (GIDEN SWIFT - Mobil) update.
$ws.Range("F14").Value = ""
